$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# =========================================================================
# 1. Shrink Tabela1 (Product Backlog) from 4 to 3 columns: drop "Coluna5"
# =========================================================================
$tabela1 = $ws1.ListObjects.Item("Tabela1")
$tabela1.ListColumns.Item(4).Delete()

# Re-style the freed column E to match the plain background column (copy from column F)
$ws1.Range("F4").Copy() | Out-Null
$ws1.Range("E4").PasteSpecial(-4122) | Out-Null
$ws1.Range("F5").Copy() | Out-Null
$ws1.Range("E5:E10").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# =========================================================================
# 2. Extend the highlighted "DoD" panel background down to row 21 and shift
#    the DoD list up by one row (header moves from D18 to D17).
# =========================================================================
$ws1.Range("A3").Copy() | Out-Null
$ws1.Range("A16:F21").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Restore the s=2 (header) look on B17, and give D17 the same look + "DoD" text
$ws1.Range("B4").Copy() | Out-Null
$ws1.Range("B17").PasteSpecial(-4122) | Out-Null
$ws1.Range("D17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws1.Range("D17").Value = "DoD"

# Restore the s=1 (white block) look on D18:D20 with the shifted DoD items
$ws1.Range("F19").Copy() | Out-Null
$ws1.Range("D18:D20").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws1.Range("D18").Value = "Codificado"
$ws1.Range("D19").Value = "Testado"
$ws1.Range("D20").Value = "Documentado"

Write-Host "Step 1-2 done"
